$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("service_tables")
$ws.Range("A1").Value = "test"
